$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill row 3 with the new strings (A3 keeps its existing style, B3 keeps its
# existing quote-prefixed style -- a leading apostrophe preserves that style
# while not becoming part of the stored text).
$ws.Range("A3").Value = "verify that user able to access the kids successfully."
$ws.Range("B3").Value = "'princess"

# New column width for column B (target raw OOXML width is 26.90625 chars;
# the host's ColumnWidth setter quantizes to 1/6-character pixel steps, so
# 26 lands on the closest reachable stored width, 26.8333...).
$ws.Columns.Item(2).ColumnWidth = 26

# Row 4: A4 keeps the same style used throughout column A; B4 gets a new
# left+vcenter aligned style built off B3's font.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A4").Value = $null

# Base B4's style on B2 (same font, vertical-center, no quote-prefix) then
# add the extra horizontal-left alignment that the new style introduces.
$ws.Range("B2").Copy()
$ws.Range("B4").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B4").HorizontalAlignment = -4131 # xlLeft
$ws.Range("B4").VerticalAlignment = -4108  # xlVAlignCenter
$ws.Range("B4").Value = $null

$ws.Range("B4").Select()
